$d = $word.ActiveDocument

# Locate the end of the existing "No such user: NotExistingUser" run; the new
# content must be inserted right after it (and before the run that currently
# holds the four spaces preceding "demonstration").
$anchor = $d.Content
$anchor.Find.Execute("No such user: NotExistingUser", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPos = $anchor.End

# Insert "<---" right after "NotExistingUser". It naturally inherits the
# red / size-32 / light-gray-highlight formatting of the text right before
# it; we only need to recolor it orange.
$arrow = "<---"
$arrowTarget = $d.Range($insertPos, $insertPos)
$arrowTarget.InsertAfter($arrow)
$arrowRange = $d.Range($insertPos, $insertPos + $arrow.Length)
$arrowRange.Font.Color = 42495

# Insert the version-mismatch message right after "<---" with the same
# orange / size-32 / light-gray-highlight formatting.
$afterArrow = $insertPos + $arrow.Length
$message = "M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0"
$messageTarget = $d.Range($afterArrow, $afterArrow)
$messageTarget.InsertAfter($message)
$messageRange = $d.Range($afterArrow, $afterArrow + $message.Length)
$messageRange.Font.Color = 42495

# Grab a formatting template from an existing, unformatted run (the four
# spaces already sitting right before "demonstration") so the new plain run
# of spaces is inserted with no run properties at all, instead of inheriting
# the orange/size/highlight formatting of the text right before it.
$template = $d.Content
$template.Find.Execute("    demonstration", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$plainTemplateStart = $template.Start
$plainTemplateRange = $d.Range($plainTemplateStart, $plainTemplateStart + 4)
$plainFormatting = $plainTemplateRange.FormattedText

# Insert the plain (unformatted) four spaces between "NotExistingUser" and
# "<---".
$spacesTarget = $d.Range($insertPos, $insertPos)
$spacesTarget.FormattedText = $plainFormatting
